# Doing Updates for Financials
# Insert two new quarterly columns (D,E) into the URI quarterly financials
# sheet, shifting the existing quarters two columns to the right (D->F,
# E->G, ... K->M), and populate the two new columns with the newest
# quarter-end figures. A couple of historical rows were also restated as
# part of this refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert two blank columns before column D. Everything that used to
#    live in D:K automatically shifts to F:M.
# ---------------------------------------------------------------------
$ws.Range("D:E").Insert()

# ---------------------------------------------------------------------
# 2) Copy the number formatting (date style / number style) from the
#    (now shifted) old column D/E -- which live in F/G after the
#    insert -- into the freshly inserted D/E columns so the new cells
#    render the same way (date format for the header rows, thousands
#    format for the data rows, etc.)
# ---------------------------------------------------------------------
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G7:G102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Fill in the new quarter values for column D (quarter ended
#    2018-12-31) and column E (quarter ended 2018-09-30).
#    Values are keyed by row number -> @(D-value, E-value).
#    "NA" means the text value NA (same shared string used elsewhere
#    in the sheet), everything else is numeric.
# ---------------------------------------------------------------------
$newQuarterValues = @{
    7   = @("43465", "43373")
    8   = @("2306000", "2116000")
    9   = @("1308000", "1178000")
    10  = @("998000", "938000")
    12  = @("NA", "NA")
    13  = @("0", "0")
    14  = @("38000", "20000")
    15  = @("95000", "75000")
    17  = @("1743000", "1538000")
    18  = @("563000", "578000")
    20  = @("-138000", "-118000")
    21  = @("895000", "878000")
    22  = @("0", "0")
    23  = @("425000", "460000")
    24  = @("109000", "127000")
    25  = @("0", "0")
    26  = @("316000", "333000")
    27  = @("316000", "333000")
    28  = @("0", "0")
    29  = @("-6000", "NA")
    30  = @("0", "0")
    31  = @("0", "0")
    32  = @("138000", "118000")
    33  = @("310000", "333000")
    34  = @("0", "0")
    35  = @("310000", "333000")
    38  = @("43465", "43373")
    41  = @("43000", "65000")
    42  = @("0", "0")
    43  = @("1545000", "1438000")
    44  = @("109000", "104000")
    45  = @("64000", "85000")
    46  = @("1761000", "1692000")
    47  = @("0", "0")
    48  = @("10214000", "9439000")
    49  = @("6142000", "5208000")
    50  = @("0", "0")
    51  = @("0", "0")
    52  = @("16000", "15000")
    53  = @("0", "0")
    54  = @("18133000", "16354000")
    57  = @("536000", "688000")
    58  = @("903000", "896000")
    59  = @("677000", "503000")
    60  = @("2116000", "2087000")
    61  = @("10844000", "9182000")
    62  = @("1770000", "1751000")
    63  = @("0", "0")
    64  = @("0", "0")
    65  = @("0", "0")
    66  = @("14730000", "13020000")
    68  = @("0", "0")
    69  = @("0", "0")
    70  = @("0", "0")
    71  = @("0", "0")
    72  = @("4101000", "3791000")
    73  = @("0", "0")
    74  = @("0", "0")
    75  = @("0", "0")
    76  = @("3403000", "3334000")
    77  = @("0", "0")
    80  = @("43465", "43373")
    81  = @("310000", "333000")
    83  = @("470000", "418000")
    84  = @("0", "0")
    85  = @("0", "0")
    86  = @("0", "0")
    87  = @("0", "0")
    88  = @("0", "0")
    89  = @("730000", "474000")
    91  = @("-51000", "-54000")
    92  = @("0", "0")
    93  = @("0", "0")
    94  = @("-2158000", "-1388000")
    96  = @("0", "0")
    97  = @("0", "0")
    98  = @("0", "0")
    99  = @("0", "0")
    100 = @("1404000", "863000")
    101 = @("2000", "-1000")
    102 = @("-22000", "-52000")
}

foreach ($row in $newQuarterValues.Keys) {
    $vals = $newQuarterValues[$row]
    $ws.Range("D$row").Value2 = $vals[0]
    $ws.Range("E$row").Value2 = $vals[1]
}

# ---------------------------------------------------------------------
# 4) Two rows in the Cash Flow Statement were restated with revised
#    historical figures rather than a simple shift of the prior values.
# ---------------------------------------------------------------------

# Row 91 - "Capital Expenditures": D..J were all restated, only the
# three oldest quarters (now K, L, M) retain their previous values.
$ws.Range("F91").Value2 = "-47000"
$ws.Range("G91").Value2 = "-33000"
$ws.Range("H91").Value2 = "-33000"
$ws.Range("I91").Value2 = "-32000"
$ws.Range("J91").Value2 = "-33000"

# Row 94 - "Total Cash Flows From Investing Activities": D,E are new,
# F,G keep the old D,E values, H,I were restated, and J,K,L,M keep the
# previously shifted values (old H,I,J,K).
$ws.Range("H94").Value2 = "-1442000"
$ws.Range("I94").Value2 = "-558000"

$wb.Save()
